## "fixed resume download button"
##
## The resume had a leftover internal review comment ("*IMPORTANT* Ensure
## that at least 3 projects showcase skills learned during the intensive
## program.") anchored on the "Applications:" run. That note was resolved/
## no longer needed, so remove it entirely (range markers + reference +
## the comment text itself) via the Word object model, exactly like a
## human reviewer right-clicking the comment and choosing "Delete Comment".

$d = $word.ActiveDocument

# Walk backwards so deleting one comment never invalidates the index of
# the ones we still have to process.
$count = $d.Comments.Count
for ($i = $count; $i -ge 1; $i--) {
    $d.Comments($i).Delete()
}
